$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refreshed crypto data (prices + 1h volume%); also rows 44/45 swap ranking
# positions (dogwifhat moves up to 44, OKB drops to 45) as part of the refresh.
#
# Column D holds price strings that look numeric (e.g. "577.97", "1.00"):
# assigning them straight to .Value lets Excel's smart-typing convert them to
# real numbers, which would not match the original text-cell (inlineStr)
# layout. Prefixing with a leading apostrophe forces literal text entry, then
# resetting .Style to "Normal" clears the quotePrefix formatting Excel
# applies for apostrophe-entered text, so the cell ends up identical in style
# to every other (untouched) text cell in the sheet.
#
# Column B/C/E values (coin names, links, padded "  +x.xx%  " volume strings)
# are never ambiguous numbers, so a plain .Value assignment is sufficient.

# --- Column D price cells: force literal text, matching original inlineStr ---
$ws.Range("D2").Value = "'67.390.89"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Value = "'3.216.08"
$ws.Range("D3").Style = "Normal"
$ws.Range("D5").Value = "'577.97"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Value = "'182.53"
$ws.Range("D6").Style = "Normal"
$ws.Range("D8").Value = "'0.602"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Value = "'3.213.91"
$ws.Range("D9").Style = "Normal"
$ws.Range("D13").Value = "'3.777.29"
$ws.Range("D13").Style = "Normal"
$ws.Range("D16").Value = "'67.466.67"
$ws.Range("D16").Style = "Normal"
$ws.Range("D18").Value = "'3.236.31"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Value = "'5.74"
$ws.Range("D19").Style = "Normal"
$ws.Range("D21").Value = "'391.94"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Value = "'7.50"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Value = "'1.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Value = "'70.73"
$ws.Range("D24").Style = "Normal"
$ws.Range("D28").Value = "'9.51"
$ws.Range("D28").Style = "Normal"
$ws.Range("D32").Value = "'22.51"
$ws.Range("D32").Style = "Normal"
$ws.Range("D36").Value = "'160.97"
$ws.Range("D36").Style = "Normal"
$ws.Range("D38").Value = "'1.88"
$ws.Range("D38").Style = "Normal"
$ws.Range("D44").Value = "'2.44"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Value = "'40.59"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Value = "'2.603.60"
$ws.Range("D46").Style = "Normal"
$ws.Range("D48").Value = "'333.72"
$ws.Range("D48").Style = "Normal"

# --- Column B/C/E cells: plain text, safe to set directly ---
$ws.Range("E2").Value = "  -0.75%  "
$ws.Range("E3").Value = "  -1.15%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("E5").Value = "  -0.97%  "
$ws.Range("E6").Value = "  -0.68%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("E8").Value = "  +0.83%  "
$ws.Range("E9").Value = "  -1.20%  "
$ws.Range("E10").Value = "  -3.32%  "
$ws.Range("E11").Value = "  -1.53%  "
$ws.Range("E12").Value = "  -0.57%  "
$ws.Range("E13").Value = "  -1.14%  "
$ws.Range("E14").Value = "  -0.14%  "
$ws.Range("E15").Value = "  -3.22%  "
$ws.Range("E17").Value = "  -2.26%  "
$ws.Range("E18").Value = "  -0.62%  "
$ws.Range("E19").Value = "  -1.94%  "
$ws.Range("E20").Value = "  -1.46%  "
$ws.Range("E21").Value = "  +2.73%  "
$ws.Range("E22").Value = "  -2.14%  "
$ws.Range("E23").Value = "  -0.01%  "
$ws.Range("E24").Value = "  -0.77%  "
$ws.Range("E25").Value = "  -0.52%  "
$ws.Range("E26").Value = "  -2.42%  "
$ws.Range("E27").Value = "  +1.16%  "
$ws.Range("E28").Value = "  -3.34%  "
$ws.Range("E29").Value = "  -0.07%  "
$ws.Range("E30").Value = "  -1.88%  "
$ws.Range("E31").Value = "  -2.40%  "
$ws.Range("E32").Value = "  -1.58%  "
$ws.Range("E33").Value = "  -4.84%  "
$ws.Range("E34").Value = "  +0.01%  "
$ws.Range("E35").Value = "  -1.70%  "
$ws.Range("E36").Value = "  -0.83%  "
$ws.Range("E38").Value = "  +1.57%  "
$ws.Range("E39").Value = "  -1.63%  "
$ws.Range("E40").Value = "  -3.90%  "
$ws.Range("E41").Value = "  -1.57%  "
$ws.Range("E42").Value = "  -4.02%  "
$ws.Range("E43").Value = "  -0.87%  "
$ws.Range("B44").Value = "dogwifhat"
$ws.Range("C44").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("E44").Value = "  -5.78%  "
$ws.Range("B45").Value = "OKB"
$ws.Range("C45").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("E45").Value = "  -1.88%  "
$ws.Range("E46").Value = "  -1.55%  "
$ws.Range("E47").Value = "  -2.99%  "
$ws.Range("E48").Value = "  -3.44%  "
$ws.Range("E49").Value = "  -2.86%  "
$ws.Range("E50").Value = "  +0.68%  "
$ws.Range("E51").Value = "  -1.81%  "
